$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.139.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.66%  "

$ws.Range("D3").Value = "'2.380.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.95%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'303.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("D6").Value = "'97.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.99%  "

$ws.Range("E7").Value = "  +0.40%  "

$ws.Range("D9").Value = "'0.501"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.00%  "

$ws.Range("D10").Value = "'34.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.89%  "

$ws.Range("D11").Value = "'0.0789"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.92%  "

$ws.Range("E12").Value = "  +2.36%  "

$ws.Range("D13").Value = "'18.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.41%  "

$ws.Range("D14").Value = "'6.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.68%  "

$ws.Range("D15").Value = "'2.754.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.93%  "

$ws.Range("D16").Value = "'2.384.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.27%  "

$ws.Range("D17").Value = "'0.810"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.96%  "

$ws.Range("D18").Value = "'43.139.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.76%  "

$ws.Range("D19").Value = "'12.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.67%  "

$ws.Range("E20").Value = "  +6.86%  "

$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("D22").Value = "'68.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.97%  "

$ws.Range("D23").Value = "'236.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.44%  "

$ws.Range("E24").Value = "  -2.34%  "

$ws.Range("D25").Value = "'2.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.29%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").Value = "'24.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.00%  "

$ws.Range("E28").Value = "  +0.61%  "

$ws.Range("D29").Value = "'9.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.14%  "

$ws.Range("D30").Value = "'31.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.96%  "

$ws.Range("E31").Value = "  +0.04%  "

$ws.Range("E32").Value = "  +2.13%  "

$ws.Range("D33").Value = "'0.0734"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.14%  "

$ws.Range("E34").Value = "  -1.56%  "

$ws.Range("E35").Value = "  +6.76%  "

$ws.Range("D36").Value = "'4.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.03%  "

$ws.Range("E37").Value = "  +2.29%  "

$ws.Range("E38").Value = "  -0.91%  "

$ws.Range("D39").Value = "'2.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.15%  "

$ws.Range("D40").Value = "'22.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.35%  "

$ws.Range("E41").Value = "  +0.31%  "

$ws.Range("D42").Value = "'108.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -34.62%  "

$ws.Range("D43").Value = "'1.951.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.61%  "

$ws.Range("E44").Value = "  +0.69%  "

$ws.Range("E45").Value = "  +1.92%  "

$ws.Range("D46").Value = "'2.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.16%  "

$ws.Range("D47").Value = "'9.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -11.35%  "

$ws.Range("D48").Value = "'2.610.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.79%  "

$ws.Range("D49").Value = "'52.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.42%  "

$ws.Range("E50").Value = "  +1.81%  "

$ws.Range("D51").Value = "'72.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.58%  "
